# Retraining the F_Cristan XGB model: update predicted GHI/DNI/DHI values
# in the Daily summary row and the Hourly breakdown rows (hours 7-17) to
# reflect the refreshed model output.

$wb = $excel.ActiveWorkbook

# --- Daily sheet: row 2 (date 2024-02-15 summary) ---
$wsDaily = $wb.Worksheets.Item("Daily")
$wsDaily.Range("G2").Value = 2989.58
$wsDaily.Range("H2").Value = 6263.94
$wsDaily.Range("I2").Value = 738.2
$wsDaily.Range("J2").Value = 2989.5
$wsDaily.Range("K2").Value = 5952.75
$wsDaily.Range("L2").Value = 733.59

# --- Hourly sheet: hour-by-hour breakdown rows 9-19 ---
$wsHourly = $wb.Worksheets.Item("Hourly")
$wsHourly.Range("L9").Value = 16.92
$wsHourly.Range("I10").Value = 441.5
$wsHourly.Range("K10").Value = 112.44
$wsHourly.Range("L10").Value = 403.82
$wsHourly.Range("M10").Value = 45.67
$wsHourly.Range("H11").Value = 254.39
$wsHourly.Range("I11").Value = 639.47
$wsHourly.Range("J11").Value = 73.04000000000001
$wsHourly.Range("K11").Value = 254.38
$wsHourly.Range("L11").Value = 618.77
$wsHourly.Range("H12").Value = 375.22
$wsHourly.Range("I12").Value = 734.95
$wsHourly.Range("J12").Value = 86.45
$wsHourly.Range("K12").Value = 375.21
$wsHourly.Range("L12").Value = 719.22
$wsHourly.Range("M12").Value = 85.5
$wsHourly.Range("H13").Value = 456.18
$wsHourly.Range("I13").Value = 782.53
$wsHourly.Range("J13").Value = 93.86
$wsHourly.Range("K13").Value = 456.17
$wsHourly.Range("L13").Value = 760.48
$wsHourly.Range("M13").Value = 97.20999999999999
$wsHourly.Range("H14").Value = 487.55
$wsHourly.Range("I14").Value = 798.61
$wsHourly.Range("J14").Value = 96.48999999999999
$wsHourly.Range("K14").Value = 487.54
$wsHourly.Range("M14").Value = 102.19
$wsHourly.Range("H15").Value = 465.87
$wsHourly.Range("I15").Value = 787.66
$wsHourly.Range("J15").Value = 94.67
$wsHourly.Range("K15").Value = 465.86
$wsHourly.Range("L15").Value = 764.5599999999999
$wsHourly.Range("M15").Value = 98.73
$wsHourly.Range("H16").Value = 393.52
$wsHourly.Range("I16").Value = 746.71
$wsHourly.Range("K16").Value = 393.52
$wsHourly.Range("L16").Value = 729.99
$wsHourly.Range("M16").Value = 87.98999999999999
$wsHourly.Range("H17").Value = 279.01
$wsHourly.Range("I17").Value = 662.38
$wsHourly.Range("K17").Value = 279.01
$wsHourly.Range("L17").Value = 643.45
$wsHourly.Range("M17").Value = 74.14
$wsHourly.Range("I18").Value = 491.94
$wsHourly.Range("K18").Value = 138.84
$wsHourly.Range("L18").Value = 464.89
$wsHourly.Range("M18").Value = 50.67
$wsHourly.Range("I19").Value = 115.36
$wsHourly.Range("L19").Value = 57.58
